# Assign AnonID in Excel instead.
# Drop the Prefix/Digits/Example helper rows, and update the Howto
# instructions to match: AnonIDs are now typed/dragged directly into the
# AnonID column, and the script is run afterwards to check/anonymize.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-13 held "Prefix:", "Digits:", a blank spacer row and "Example:"
# (with its helper formula). Deleting them shifts the data table header
# (previously on row 14) up to row 9.
$ws.Range("A9:R13").EntireRow.Delete() | Out-Null

# Update Howto steps 2 and 3 to describe assigning AnonIDs directly in
# Excel instead of via a Prefix/Digits scheme.
$ws.Range("D4").Value = '2. Assign AnonIDs to your Persons. You can for example type in the first (eg "MYPROJ-001") and then drag down to number the rest sequentially.'
$ws.Range("D5").Value = "3. Run aida-pat-anonexcel.py <path to this file> to check for mistakes, find and anonymize slides, and update this sheet to match."

# Match the author's final selection.
$ws.Range("D5").Select() | Out-Null
